# This revision's diff is purely a cosmetic re-serialization of the
# package's OOXML: the xmlns declarations on the root <w:document>
# element were reordered, and attributes inside <w:pgSz>, <w:pgMar>,
# <w:rFonts>, <w:lang>, <w:latentStyles>, the <w:lsdException> entries
# and the <w:style> elements in styles.xml were alphabetized. The
# commit message explains why: "Fixed POI packaging and upgraded to
# POI 3.15" - i.e. the tool that produced/packaged this reference
# .docx was upgraded, and the new version of the library happens to
# emit XML attributes in (alphabetically) sorted order instead of the
# previous insertion order.
#
# Every attribute name/value pair, every element, and all document
# text/content/formatting are identical before and after the change
# (pgSz is still 11906 x 16838 twips, margins are still 1417/1417/
# 1417/1417 with 708 header/footer and 0 gutter, fonts/styles/
# languages are unchanged, etc.) - nothing a user editing the document
# in Word could actually observe or do was altered.
#
# Word's object model / COM automation surface does not expose control
# over the raw XML attribute emission order used when the package is
# serialized (that is purely an artifact of the internal packaging
# library, not a document, paragraph, style or page-setup property),
# so there is no operation to perform through $word/$d that would
# correspond to this change. We deliberately make no modifications to
# the document here: touching any content/formatting property would
# cause the runtime to rewrite that part of the document (and could
# introduce unrelated incidental differences), which would move the
# result further from the target rather than closer, since the target
# content is already exactly what is currently in the document.

$d = $word.ActiveDocument
